$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "3BDS"
$ws.Range("B9").Value = 1029384
